$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.199.85"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.601.74"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3779"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08119"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.581"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.387"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001248"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "1.603.21"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06874"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.528"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "23.206.56"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.390"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.987"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.242"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.418"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "1.779.96"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07518"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.116"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2502"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08795"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.361"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6534"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.305"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.011"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07952"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.99%  "
